$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.081.35'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '1.731.33'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.79'
$ws.Range("E5").Value = '  -5.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4863'
$ws.Range("E7").Value = '  +6.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3515'
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.94'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.055'
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.899'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '1.730.77'
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.905'
$ws.Range("E16").Value = '  -3.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.41'
$ws.Range("E17").Value = '  -4.59%  '
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06405'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9994'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.60'
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = '27.136.27'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.90'
$ws.Range("E24").Value = '  -2.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.080'
$ws.Range("E25").Value = '  -3.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.90'
$ws.Range("E26").Value = '  -5.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.99'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '1.930.13'
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.094'
$ws.Range("E29").Value = '  -2.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.67'
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09307'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.632'
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.426'
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02200'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05933'
$ws.Range("E36").Value = '  -2.74%  '
$ws.Range("E37").Value = '  -5.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.431'
$ws.Range("E38").Value = '  +4.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.795'
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("E40").Value = '  -2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6020'
$ws.Range("E41").Value = '  -2.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9992'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.096'
$ws.Range("E43").Value = '  -7.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.521'
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.90'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.586'
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5681'
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.98'
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.853'
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.112'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06652'
$ws.Range("E51").Value = '  -1.90%  '
